# Updated symbol list on Wed Jan 25 10:44:00 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) text values for
# the crypto rows that changed. Values in this sheet are stored as plain
# text (e.g. "300.78", "-4.63%"), so we force each target cell to stay text
# (NumberFormat "@") before assigning the new literal, then restore the
# cell's style to "Normal" so no stray numeric/percentage formatting is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $newValue
    $range.Style = "Normal"
}

Set-TextValue "D2" "300.78"
Set-TextValue "E2" "-4.63%"
Set-TextValue "D3" "35.17"
Set-TextValue "E3" "-1.37%"
Set-TextValue "D4" "5.046"
Set-TextValue "E4" "-1.40%"
Set-TextValue "D5" "0.07945"
Set-TextValue "E5" "-2.01%"
Set-TextValue "D6" "1.902"
Set-TextValue "E6" "-11.05%"
Set-TextValue "D7" "7.794"
Set-TextValue "E7" "-2.62%"
Set-TextValue "D8" "4.036"
Set-TextValue "E8" "-2.73%"
Set-TextValue "E9" "4.13%"
Set-TextValue "D10" "0.9208"
Set-TextValue "E10" "-0.52%"
Set-TextValue "D11" "0.1324"
Set-TextValue "E11" "29.01%"
Set-TextValue "D12" "0.1844"
Set-TextValue "E12" "-1.74%"
Set-TextValue "D13" "0.09499"
Set-TextValue "E13" "3.26%"
Set-TextValue "D14" "0.03604"
Set-TextValue "E14" "0.44%"
Set-TextValue "D15" "0.09846"
Set-TextValue "E15" "-0.63%"
Set-TextValue "D16" "0.001395"
Set-TextValue "E16" "-2.50%"
Set-TextValue "D17" "0.005751"
Set-TextValue "E17" "1.26%"
Set-TextValue "E18" "0.97%"
Set-TextValue "D19" "0.3426"
Set-TextValue "E19" "1.75%"
Set-TextValue "E20" "-1.49%"
Set-TextValue "D21" "5.043"
Set-TextValue "E21" "-1.68%"
Set-TextValue "E22" "11.15%"
Set-TextValue "D23" "0.04495"
Set-TextValue "E23" "-1.69%"
Set-TextValue "D24" "0.001220"
Set-TextValue "E24" "-2.14%"
Set-TextValue "E25" "1.68%"
Set-TextValue "E26" "0.09%"
Set-TextValue "D27" "0.0003006"
Set-TextValue "E27" "-33.26%"
Set-TextValue "D39" "0.01872"
Set-TextValue "E39" "-4.62%"
Set-TextValue "D40" "0.04711"
Set-TextValue "E40" "-3.11%"
Set-TextValue "D41" "0.007540"
Set-TextValue "E41" "-2.28%"
Set-TextValue "D42" "0.009734"
Set-TextValue "E42" "24.33%"
Set-TextValue "D43" "0.1323"
Set-TextValue "E43" "-4.84%"
Set-TextValue "E44" "-1.77%"
Set-TextValue "D45" "0.009611"
Set-TextValue "E45" "-17.41%"
Set-TextValue "D46" "0.00006209"
Set-TextValue "E46" "-4.70%"
Set-TextValue "E47" "0.09%"
Set-TextValue "E48" "75.94%"
Set-TextValue "E49" "-12.46%"
Set-TextValue "D50" "0.00002104"
Set-TextValue "E50" "0.09%"
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.09%"
